$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 727-728; existing rows 727+ shift down to 729+.
$ws.Rows("727:728").Insert()

# Populate new row 727
$ws.Range("A727").Value = 5
$ws.Range("B727").Value = "Macroferia Regional de Talca"
$ws.Range("C727").Value = "Maule"
$ws.Range("D727").Value = 44918
$ws.Range("E727").Value = 7
$ws.Range("F727").Value = 100112020
$ws.Range("G727").Value = "Tomate"
$ws.Range("H727").Value = "Larga vida"
$ws.Range("I727").Value = "Primera"
$ws.Range("J727").Value = 2500
$ws.Range("K727").Value = 13000
$ws.Range("L727").Value = 13000
$ws.Range("M727").Value = 13000
$ws.Range("N727").Value = "`$/bandeja 18 kilos"
$ws.Range("O727").Value = "Región del Maule"
$ws.Range("P727").Value = 722
$ws.Range("Q727").Value = 18
$ws.Range("R727").Value = "Hortaliza"

# Populate new row 728
$ws.Range("A728").Value = 5
$ws.Range("B728").Value = "Macroferia Regional de Talca"
$ws.Range("C728").Value = "Maule"
$ws.Range("D728").Value = 44918
$ws.Range("E728").Value = 7
$ws.Range("F728").Value = 100112020
$ws.Range("G728").Value = "Tomate"
$ws.Range("H728").Value = "Larga vida"
$ws.Range("I728").Value = "Primera"
$ws.Range("J728").Value = 2500
$ws.Range("K728").Value = 7000
$ws.Range("L728").Value = 7000
$ws.Range("M728").Value = 7000
$ws.Range("N728").Value = "`$/caja 12 kilos"
$ws.Range("O728").Value = "Región del Maule"
$ws.Range("P728").Value = 583
$ws.Range("Q728").Value = 12
$ws.Range("R728").Value = "Hortaliza"
